$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("assuntos")

$ws.Range("A10").Value = "Conjunção"
$ws.Range("B10").Value = "Português"
$ws.Range("C10").Value = 9
